$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns stay plain text (values like "2.45" or
# "0.0618" would otherwise be auto-converted to numbers by Excel).
$ws.Range("D2:E51").NumberFormat = "@"

# --- Row swap: Chainlink (row22) <-> Toncoin (row23) ---
$ws.Range("B22").Value = "Toncoin"
$ws.Range("C22").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D22").Value = "2.32"
$ws.Range("E22").Value = "  -8.81%  "
$ws.Range("B23").Value = "Chainlink"
$ws.Range("C23").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D23").Value = "6.14"
$ws.Range("E23").Value = "  -1.90%  "

# --- Row swap: LidoDAOToken (row35) <-> HuobiToken (row36) ---
$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").Value = "2.45"
$ws.Range("E35").Value = "  +0.03%  "
$ws.Range("B36").Value = "LidoDAOToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D36").Value = "1.51"
$ws.Range("E36").Value = "  -2.44%  "

# --- Remaining price / volume updates ---
$ws.Range("D2").Value = "26.709.08"
$ws.Range("E2").Value = "  -0.21%  "
$ws.Range("D3").Value = "1.635.14"
$ws.Range("E3").Value = "  -0.82%  "
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").Value = "217.62"
$ws.Range("E5").Value = "  +0.49%  "
$ws.Range("D6").Value = "0.497"
$ws.Range("E6").Value = "  -1.52%  "
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("E8").Value = "  -1.48%  "
$ws.Range("D9").Value = "0.0618"
$ws.Range("E9").Value = "  -1.27%  "
$ws.Range("E10").Value = "  -1.59%  "
$ws.Range("D11").Value = "0.0842"
$ws.Range("E11").Value = "  -0.04%  "
$ws.Range("D12").Value = "1.867.25"
$ws.Range("E12").Value = "  -0.62%  "
$ws.Range("D13").Value = "1.639.19"
$ws.Range("E13").Value = "  -0.50%  "
$ws.Range("D14").Value = "4.10"
$ws.Range("E14").Value = "  -2.08%  "
$ws.Range("D15").Value = "0.521"
$ws.Range("E15").Value = "  -2.28%  "
$ws.Range("D16").Value = "63.89"
$ws.Range("E16").Value = "  -2.15%  "
$ws.Range("D17").Value = "26.675.91"
$ws.Range("E17").Value = "  -0.40%  "
$ws.Range("D18").Value = "0.0₃0719"
$ws.Range("E18").Value = "  -3.33%  "
$ws.Range("E19").Value = "  +0.13%  "
$ws.Range("D20").Value = "208.41"
$ws.Range("E20").Value = "  -4.17%  "
$ws.Range("E21").Value = "  -1.39%  "
$ws.Range("D24").Value = "9.15"
$ws.Range("E24").Value = "  -3.49%  "
$ws.Range("D25").Value = "146.95"
$ws.Range("E25").Value = "  +0.16%  "
$ws.Range("E26").Value = "  +0.11%  "
$ws.Range("E27").Value = "  -2.85%  "
$ws.Range("D28").Value = "6.99"
$ws.Range("E28").Value = "  -2.66%  "
$ws.Range("D29").Value = "15.44"
$ws.Range("E29").Value = "  -2.10%  "
$ws.Range("D30").Value = "0.0497"
$ws.Range("E30").Value = "  -4.27%  "
$ws.Range("E31").Value = "  +0.05%  "
$ws.Range("E32").Value = "  -1.17%  "
$ws.Range("D33").Value = "2.93"
$ws.Range("E33").Value = "  -2.64%  "
$ws.Range("D34").Value = "1.260.19"
$ws.Range("E34").Value = "  -1.57%  "
$ws.Range("D37").Value = "0.0172"
$ws.Range("E37").Value = "  -3.79%  "
$ws.Range("D38").Value = "0.518"
$ws.Range("E38").Value = "  -3.98%  "
$ws.Range("E39").Value = "  +0.05%  "
$ws.Range("D40").Value = "0.794"
$ws.Range("E40").Value = "  -4.23%  "
$ws.Range("E41").Value = "  -1.89%  "
$ws.Range("E42").Value = "  -3.20%  "
$ws.Range("D43").Value = "1.778.13"
$ws.Range("E43").Value = "  -0.73%  "
$ws.Range("D44").Value = "5.23"
$ws.Range("D45").Value = "90.82"
$ws.Range("E45").Value = "  -1.31%  "
$ws.Range("D46").Value = "59.49"
$ws.Range("E46").Value = "  -0.36%  "
$ws.Range("D47").Value = "1.56"
$ws.Range("E47").Value = "  -2.78%  "
$ws.Range("E48").Value = "  -2.06%  "
$ws.Range("D49").Value = "0.0518"
$ws.Range("E49").Value = "  +0.51%  "
$ws.Range("E50").Value = "  +0.12%  "
$ws.Range("E51").Value = "  -0.18%  "
